$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 974; this shifts rows 974:1011 down to 975:1012
$ws.Rows("974:974").Insert()

# Populate the newly inserted row 974 with the new weekly record.
# Columns A-T: Mercado ID, Mercado, Region, Fecha, Codreg, Tipo, Producto ID,
# Producto, Categoria ID, Categoria, Variedad, Calidad, Volumen, Precio minimo,
# Precio maximo, Precio promedio ponderado, Unidad de comercializacion, Origen,
# Precio $/Kg, Kg / unidad
$ws.Cells.Item(974, 1).Value = 10
$ws.Cells.Item(974, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(974, 3).Value = "La Araucanía"
$ws.Cells.Item(974, 4).Value = 45147
$ws.Cells.Item(974, 5).Value = 9
$ws.Cells.Item(974, 6).Value = "Fruta"
$ws.Cells.Item(974, 7).Value = 100102
$ws.Cells.Item(974, 8).Value = "Cítricos"
$ws.Cells.Item(974, 9).Value = 100102004
$ws.Cells.Item(974, 10).Value = "Mandarina"
$ws.Cells.Item(974, 11).Value = "Murcott"
$ws.Cells.Item(974, 12).Value = "Primera"
$ws.Cells.Item(974, 13).Value = 100
$ws.Cells.Item(974, 14).Value = 14000
$ws.Cells.Item(974, 15).Value = 14000
$ws.Cells.Item(974, 16).Value = 14000
$ws.Cells.Item(974, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(974, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(974, 19).Value = 778
$ws.Cells.Item(974, 20).Value = 18

# Ensure the date cell keeps the same date number format as the rest of column D
$ws.Cells.Item(974, 4).NumberFormat = $ws.Cells.Item(975, 4).NumberFormat
